$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
# Row 19
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()

# Row 55
$ws.Range("H55").Value = 694.1667
$ws.Range("I55").Value = 844
$ws.Range("K55").Value = 844
$ws.Range("M55").Value = -630

# Row 112
$ws.Range("H112").Value = 2122.6667
$ws.Range("I112").Value = 1500
$ws.Range("J112").Value = 2153.8
$ws.Range("K112").Value = 4500
$ws.Range("L112").Value = 6461.400000000001
$ws.Range("M112").Value = -3392
$ws.Range("N112").Value = -8677.400000000001

# Row 132
$ws.Range("H132").Value = 4222.7437
$ws.Range("I132").Value = 4202.2896
$ws.Range("K132").Value = 12606.8688
$ws.Range("M132").Value = -10076.8688

# Row 138
$ws.Range("H138").Value = 3176.309
$ws.Range("I138").Value = 4087.0557
$ws.Range("J138").Value = 2733.2432
$ws.Range("K138").Value = 12261.1671
$ws.Range("L138").Value = 8199.729599999999
$ws.Range("M138").Value = -7121.167099999999
$ws.Range("N138").Value = -18479.7296

# Row 141
$ws.Range("H141").Value = 2826.0833
$ws.Range("I141").Value = 2826.0833
$ws.Range("K141").Value = 8478.249899999999
$ws.Range("M141").Value = -3298.249899999999


# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
# Row 5
$ws.Range("H5").Value = 129.1
$ws.Range("I5").Value = 129.1
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 129.1
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -17.09999999999999
$ws.Range("N5").ClearContents()

# Row 32
$ws.Range("H32").Value = 5656.4375
$ws.Range("I32").Value = 5293
$ws.Range("K32").Value = 5293
$ws.Range("M32").Value = -5006

# Row 61
$ws.Range("H61").Value = 3034.625
$ws.Range("I61").Value = 1997
$ws.Range("K61").Value = 1997
$ws.Range("M61").Value = -1785

# Row 102
$ws.Range("H102").Value = 5714.4
$ws.Range("I102").Value = 5627.1113
$ws.Range("K102").Value = 5627.1113
$ws.Range("M102").Value = -4005.1113

# Row 132
$ws.Range("H132").Value = 2802.6296
$ws.Range("I132").Value = 1729.7646
$ws.Range("K132").Value = 5189.293799999999
$ws.Range("M132").Value = -2659.293799999999

# Row 136
$ws.Range("H136").Value = 3034.625
$ws.Range("I136").Value = 1997
$ws.Range("K136").Value = 5991
$ws.Range("M136").Value = -3441


# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
# Row 4
$ws.Range("H4").Value = 129.1
$ws.Range("I4").Value = 129.1
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 129.1
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -14.09999999999999
$ws.Range("N4").ClearContents()

# Row 105
$ws.Range("H105").Value = 17335324
$ws.Range("J105").Value = 41668916
$ws.Range("L105").Value = 41668916
$ws.Range("N105").Value = -41672410

# Row 134
$ws.Range("H134").Value = 3103.5
$ws.Range("I134").Value = 2811.756
$ws.Range("K134").Value = 8435.268
$ws.Range("M134").Value = -5900.268


# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
# Row 31
$ws.Range("H31").Value = 4853.909
$ws.Range("I31").Value = 3020.5715
$ws.Range("K31").Value = 3020.5715
$ws.Range("M31").Value = -2725.5715

# Row 34
$ws.Range("H34").Value = 4853.909
$ws.Range("I34").Value = 3020.5715
$ws.Range("K34").Value = 3020.5715
$ws.Range("M34").Value = -2818.5715

# Row 55
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()


# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
# Row 5
$ws.Range("H5").Value = 1095.4117
$ws.Range("I5").Value = 785
$ws.Range("J5").Value = 1371.3334
$ws.Range("K5").Value = 2355
$ws.Range("L5").Value = 4114.0002
$ws.Range("M5").Value = -2243
$ws.Range("N5").Value = -4338.0002

# Row 12
$ws.Range("H12").Value = 124.9
$ws.Range("J12").Value = 97.666664
$ws.Range("L12").Value = 292.999992
$ws.Range("N12").Value = -638.999992

# Row 36
$ws.Range("H36").Value = 1899.75
$ws.Range("J36").Value = 3349.5
$ws.Range("L36").Value = 10048.5
$ws.Range("N36").Value = -10386.5

# Row 122
$ws.Range("H122").Value = 561.13635
$ws.Range("J122").Value = 451.6
$ws.Range("L122").Value = 4064.4
$ws.Range("N122").Value = -8964.4

# Row 132
$ws.Range("H132").Value = 1869.579
$ws.Range("I132").Value = 1450.4286
$ws.Range("J132").Value = 2114.0833
$ws.Range("K132").Value = 13053.8574
$ws.Range("L132").Value = 19026.7497
$ws.Range("M132").Value = -10523.8574
$ws.Range("N132").Value = -24086.7497

# Row 135
$ws.Range("H135").Value = 1095.4117
$ws.Range("I135").Value = 785
$ws.Range("J135").Value = 1371.3334
$ws.Range("K135").Value = 7065
$ws.Range("L135").Value = 12342.0006
$ws.Range("M135").Value = -4530
$ws.Range("N135").Value = -17412.0006

# Row 138
$ws.Range("H138").Value = 1128.75
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()


# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
# Row 42
$ws.Range("H42").Value = 49103
$ws.Range("J42").Value = 49103
$ws.Range("L42").Value = 49103
$ws.Range("N42").Value = -50073

# Row 70
$ws.Range("H70").Value = 85505.24000000001
$ws.Range("I70").Value = 170426.83
$ws.Range("J70").Value = 7116.077
$ws.Range("K70").Value = 170426.83
$ws.Range("L70").Value = 7116.077
$ws.Range("M70").Value = -170156.83
$ws.Range("N70").Value = -7656.077

# Row 73
$ws.Range("H73").Value = 85505.24000000001
$ws.Range("I73").Value = 170426.83
$ws.Range("J73").Value = 7116.077
$ws.Range("K73").Value = 170426.83
$ws.Range("L73").Value = 7116.077
$ws.Range("M73").Value = -169490.83
$ws.Range("N73").Value = -8988.077000000001

# Row 115
$ws.Range("H115").Value = 49103
$ws.Range("J115").Value = 49103
$ws.Range("L115").Value = 49103
$ws.Range("N115").Value = -51453


# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
# Row 46
$ws.Range("H46").Value = 3383.3333
$ws.Range("J46").Value = 4750
$ws.Range("L46").Value = 4750
$ws.Range("N46").Value = -5126

# Row 61
$ws.Range("H61").Value = 7097.5264
$ws.Range("I61").Value = 1685.0834
$ws.Range("J61").Value = 16376
$ws.Range("K61").Value = 1685.0834
$ws.Range("L61").Value = 16376
$ws.Range("M61").Value = -1483.0834
$ws.Range("N61").Value = -16780

# Row 113
$ws.Range("H113").Value = 7097.5264
$ws.Range("I113").Value = 1685.0834
$ws.Range("J113").Value = 16376
$ws.Range("K113").Value = 1685.0834
$ws.Range("L113").Value = 16376
$ws.Range("M113").Value = 484.9166
$ws.Range("N113").Value = -20716

# Row 122
$ws.Range("H122").Value = 9897.714
$ws.Range("I122").Value = 9219.556
$ws.Range("J122").Value = 11118.4
$ws.Range("K122").Value = 27658.668
$ws.Range("L122").Value = 33355.2
$ws.Range("M122").Value = -25208.668
$ws.Range("N122").Value = -38255.2

# Row 132
$ws.Range("H132").Value = 4632.346
$ws.Range("I132").Value = 3473.647
$ws.Range("J132").Value = 6821
$ws.Range("K132").Value = 10420.941
$ws.Range("L132").Value = 20463
$ws.Range("M132").Value = -7890.940999999999
$ws.Range("N132").Value = -25523

# Row 136
$ws.Range("H136").Value = 4792.5713
$ws.Range("I136").Value = 5912.125
$ws.Range("K136").Value = 17736.375
$ws.Range("M136").Value = -15186.375

# Row 138
$ws.Range("H138").Value = 115927.5
$ws.Range("J138").Value = 115927.5
$ws.Range("L138").Value = 115927.5
$ws.Range("N138").Value = -126207.5

# Row 141
$ws.Range("H141").Value = 99999
$ws.Range("J141").Value = 99999
$ws.Range("L141").Value = 99999
$ws.Range("N141").Value = -110359


# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
# Row 107
$ws.Range("H107").Value = 833.8570999999999
$ws.Range("I107").Value = 856.25
$ws.Range("K107").Value = 2568.75
$ws.Range("M107").Value = -648.75

# Row 113
$ws.Range("H113").Value = 418.26666
$ws.Range("I113").Value = 241.9
$ws.Range("K113").Value = 725.7
$ws.Range("M113").Value = 1444.3

# Row 129
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

# Row 132
$ws.Range("H132").Value = 6746
$ws.Range("I132").Value = 7449.6665
$ws.Range("J132").Value = 6142.857
$ws.Range("K132").Value = 22348.9995
$ws.Range("L132").Value = 18428.571
$ws.Range("M132").Value = -19818.9995
$ws.Range("N132").Value = -23488.571

